$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = "BSc FT Computer Science"
$ws.Range("K2").Value = "School of X"
$ws.Range("S2").Value = "Female"
$ws.Range("AD2").Value = "Yes"
$ws.Range("AE2").Value = "Speaks Urdu, interested in programming"
$ws.Range("AG2").Value = "Software Engineering or Programming, Cyber Security"
$ws.Range("AN2").Value = "Planning for the future and goal setting, Gaining insight to an industry/profession, Building a professional network, Writing/improving CVs, job applications and covering letters"
$ws.Range("AO2").Value = "No Preference"
$ws.Range("AP2").Value = "Option 1 - A mentor who studied the same degree as me but works in any industry/job role"

# Row 3
$ws.Range("I3").Value = "BA FT Drama"
$ws.Range("K3").Value = "School of Y"
$ws.Range("S3").Value = "Male"
$ws.Range("AD3").Value = "No"
$ws.Range("AE3").Value = "none"
$ws.Range("AG3").Value = "English or writer"
$ws.Range("AN3").Value = "Planning for the future and goal setting, Gaining insight to an industry/profession, Building a professional network, Writing/improving CVs, job applications and covering letters"
$ws.Range("AO3").Value = "Female"
$ws.Range("AP3").Value = "Option 2 - A mentor who works in the industry/job role that I am interested in"

# Row 4
$ws.Range("I4").Value = "BSc FT Chemistry"
$ws.Range("K4").Value = "School of X"
$ws.Range("S4").Value = "Female"
$ws.Range("AD4").Value = "No"
$ws.Range("AE4").Value = "no"
$ws.Range("AG4").Value = "chemist, pharmasist"
$ws.Range("AI4").Value = "Stage 1 - I am looking to explore entrepreneurship"
$ws.Range("AN4").Value = "Interview practice and preparation, Finding work experience (shadowing/internships/part-time work), Developing entrepreneurial skills, Support with setting up or growing a business"
$ws.Range("AO4").Value = "Male"
$ws.Range("AP4").Value = "Option 3 - A mentor who can support with entrepreneurship"

# Row 5 (new row)
$ws.Range("D5").Value = "A"
$ws.Range("E5").Value = "C"
$ws.Range("G5").Value = "13242R455"
$ws.Range("I5").Value = "BSc FT Chemistry"
$ws.Range("K5").Value = "School of S"
$ws.Range("S5").Value = "Female"
$ws.Range("AD5").Value = "No"
$ws.Range("AG5").Value = "scientist, researcher, biologist"
$ws.Range("AI5").Value = "Stage 3 - I am a current entrepreneur in need extra support with my business"
$ws.Range("AN5").Value = "Support with setting up or growing a business, Planning for the future and goal setting, Building a professional network"
$ws.Range("AO5").Value = "Female"
$ws.Range("AP5").Value = "Option 3 - A mentor who can support with entrepreneurship"

# Sheet view changes: scroll the window so column AM is at the left edge,
# zoom to 119%, and leave AM4 as the active/selected cell.
$ws.Range("AM4").Select()
$win = $ws.Application.ActiveWindow
$win.ScrollColumn = 39
$win.ScrollRow = 1
$win.Zoom = 119
